$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e0c4adeb60343a2fd2a25a9a6bc0749afa108b/e2e/6a062ce1-3fdf-43db-bb35-23771fe24cbb.md"
$display = "6a062ce1-3fdf-43db-bb35-23771fe24cbb.md"
$errMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cead5ca0bdc81c0f66255625561c7dfe38bf0110/e2e/6a062ce1-3fdf-43db-bb35-23771fe24cbb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a4e0c4adeb60343a2fd2a25a9a6bc0749afa108b/e2e/6a062ce1-3fdf-43db-bb35-23771fe24cbb.md."

# --- zh-cn sheet, row 7 (6a062ce1-... entry) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, [Type]::Missing, [Type]::Missing, $display)
$wsZh.Range("J7").Value2 = "6a062ce1-3fdf-43db-bb35-23771fe24cbb.e800fd154e2d2aa33ac1b1a0d16f46ce6ce69ee6.zh-cn.xlf"
$wsZh.Range("K7").Value2 = "2016-09-01 11:04:31"
$wsZh.Range("P7").Value2 = $errMsg

# --- de-de sheet, row 7 (6a062ce1-... entry) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, [Type]::Missing, [Type]::Missing, $display)
$wsDe.Range("J7").Value2 = "6a062ce1-3fdf-43db-bb35-23771fe24cbb.e800fd154e2d2aa33ac1b1a0d16f46ce6ce69ee6.de-de.xlf"
$wsDe.Range("K7").Value2 = "2016-09-01 11:04:38"
$wsDe.Range("P7").Value2 = $errMsg
